# Update automatico via Actualizar 02-06-2021 05-40-06
#
# Appends one new "availability check" block (14 rows, one per monitored
# service) to Sheet1, and refreshes the timestamp of the previous block
# (rows 1024-1037) to reflect the precise recalculated serial time.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-stamp the previous block (rows 1024-1037) with the refined time ---
$refreshedSerial = 44233.21495916667
for ($r = 1024; $r -le 1037; $r++) {
    $ws.Range("D$r").Value = $refreshedSerial
}

# --- 2. Append the new block (rows 1038-1051) ---
$newSerial = 44233.23612935856

$entries = @(
    @{ Name = "Odoo";              Url = "https://www.dataintelligence-group.com/" },
    @{ Name = "Blackbox";          Url = "https://serviciodashboard.azurewebsites.net/" },
    @{ Name = "PowerBI";           Url = "https://powerbi.microsoft.com/es-es/" },
    @{ Name = "Dropbox";           Url = "https://www.dropbox.com/" },
    @{ Name = "Odoo";              Url = "https://dataintelligence.store/" },
    @{ Name = "GEE";               Url = "https://app-data-i.users.earthengine.app/" },
    @{ Name = "UtilidadesOdoo";    Url = "https://odooutil.azurewebsites.net/" },
    @{ Name = "Filtros Dashboard"; Url = "https://filtradordashboard.azurewebsites.net/" },
    @{ Name = "MapStore";          Url = "https://ide.dataintelligence-group.com/mapstore/#/" },
    @{ Name = "GeoServer";         Url = "https://ide.dataintelligence-group.com/geoserver/web/?0" },
    @{ Name = "Tomcat";            Url = "https://ide.dataintelligence-group.com/" },
    @{ Name = "Shiny";             Url = "https://rpubs.com/dataintelligence/" },
    @{ Name = "Github";            Url = "https://github.com/Sud-Austral/" },
    @{ Name = "EZ Exporter";       Url = "https://ezexporter.highviewapps.com/exports/export-profile/" }
)

$row = 1038
foreach ($entry in $entries) {
    $ws.Range("A$row").Value = $entry.Name
    $ws.Range("B$row").Value = $entry.Url

    # MapStore's URL carries a "#/" fragment, which Excel treats as the
    # hyperlink's sub-address (location within the page) rather than part
    # of the target, matching the pattern already used for earlier blocks.
    $hashIndex = $entry.Url.IndexOf("#")
    if ($hashIndex -ge 0) {
        $address = $entry.Url.Substring(0, $hashIndex)
        $subAddress = $entry.Url.Substring($hashIndex + 1)
        $ws.Hyperlinks.Add($ws.Range("B$row"), $address, $subAddress)
    } else {
        $ws.Hyperlinks.Add($ws.Range("B$row"), $entry.Url)
    }
    $ws.Range("B$row").Style = "Hyperlink"

    $ws.Range("C$row").Value = "Disponible"
    $ws.Range("D$row").Value = $newSerial
    $ws.Range("D$row").NumberFormat = "YYYY-MM-DD HH:MM:SS"

    $row++
}
